$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix device connection issues: update hostnames/IPs for sandbox devices
$ws.Range("A2").Value = "sbx-ao"
$ws.Range("B3").Value = "ios-xe-mgmt.cisco.com"

# Update the saved selection on the sheet to B3
$ws.Range("B3").Select()
